# "Generate Report for Archive"
# The localization status for a0df9861-1dd6-4eb5-a6c4-8ed409b650af.md has
# moved on from handoff: it is now back "In Translation" (instead of
# "Ready for handoff") for every locale, both on the per-locale report
# sheets and on the roll-up Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"
